# Display BARS Code in annual report.
#
# The "Details" sheet gets a new "BARS Code" column inserted right after
# "Account Description" (i.e. a brand-new column C, pushing the existing
# Act #1..#5 / Debit / Credit / Amount columns one slot to the right).
# The "Details" tab also becomes the active/selected sheet, and the
# previously-active "Schedule 1" sheet keeps a remembered selection of C5.

$wb = $excel.ActiveWorkbook

$schedule1 = $wb.Worksheets.Item("Schedule 1")
$details   = $wb.Worksheets.Item("Details")

# Insert a new column before column C ("Act #1") on the Details sheet and
# give it the "BARS Code" header to match the rest of row 1's formatting.
$details.Range("C1").EntireColumn.Insert()
$details.Range("C1").Value = "BARS Code"

# Remember a selection on "Schedule 1" (the sheet that used to be active)
# before switching the active tab over to "Details".
$schedule1.Range("C5").Select()
$details.Activate()
